$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 24, pushing existing rows 24-26 down to 25-27
$ws.Rows("24:24").Insert()

# Populate the new row 24 with the new data record
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44754
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112042
$ws.Range("G24").Value = "Locoto"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 3300
$ws.Range("L24").Value = 3300
$ws.Range("M24").Value = 3300
$ws.Range("N24").Value = "$/kilo"
$ws.Range("O24").Value = "Región de Arica y Parinacota"
$ws.Range("P24").Value = 3300
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"
